# Integrate retry logic and data provider enhancements; add method
# interceptor and annotation transformer for test execution control.
#
# This adds a new "TestsRunner" worksheet (after the existing "Sheet1")
# that drives which TestNG test cases should execute, with a short
# description for each.

$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the workbook (after Sheet1) and name it.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "TestsRunner"

# Populate header + data. Cells are written in this specific order so the
# workbook's shared-string table is built up the same way a person typing
# the sheet column-by-column, row-by-row (but filling the Execute flags as
# they went) would produce it.
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Description"

$ws.Range("A2").Value = "loginLogoutTest"
$ws.Range("B2").Value = "validate OrangeHRM login and logout functionality"
$ws.Range("C3").Value = "yes"

$ws.Range("A3").Value = "homePageTitleTest"
$ws.Range("C2").Value = "no"

$ws.Range("B3").Value = "validate title of home page"
$ws.Range("C1").Value = "Execute"

# Widen the TestCase / Description columns so the text isn't clipped.
$ws.Columns.Item(1).ColumnWidth = 20.14
$ws.Columns.Item(2).ColumnWidth = 47.46

# Leave the cursor on C4, as it was when the sheet was last saved.
$ws.Range("C4").Select() | Out-Null

Write-Output "TestsRunner sheet added"
